$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh crypto price (column D) and volume/1h change (column E) figures.
# Force text formatting on these cells so values stay exact strings
# (matching the source data's inline-string representation) rather than
# being auto-coerced into numeric cells by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.409.66"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.57%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.144.89"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.90%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "526.45"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -4.92%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.99"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.40%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.146.19"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.77%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.443"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -4.33%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.20"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -7.22%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -7.55%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.380"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -5.84%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.679.21"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.38%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.92%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.58"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.43%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.141.56"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -4.49%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "58.377.56"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.81%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000153"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -6.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.78"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.05"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -5.04%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.93"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -6.75%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "344.93"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -6.89%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.508"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.82%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "68.00"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -7.59%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.266.07"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.87%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.169"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.33%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0963"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.95%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.996"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.82"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.53%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.87"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -7.35%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.88"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -7.35%  "
# Rows 34 and 35 swap identity: Fetch.AI and EthereumClassic trade places
# in the ranking, and each coin's price/volume figures are refreshed.
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.23"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.49%  "

$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "21.43"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.31%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.82"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.18%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "157.26"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.26"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -5.19%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.38"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -8.49%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0688"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.87%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.172.04"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.15%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -5.85%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.46"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.71%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.84%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.693"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -6.79%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.91"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.90%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.02%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -7.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.273.89"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.27%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.22"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.58%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.81"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.35%  "
